$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '62.043.03'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = '2.440.11'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  -0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '579.62'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +2.02%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '142.95'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -1.53%  '
$ws.Range('E7').Value = '  +0.00%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.529'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -0.65%  '
$ws.Range('D9').Value = '2.438.24'
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('E10').Value = '  -2.98%  '
$ws.Range('E11').Value = '  +2.38%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '5.18'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -1.00%  '
$ws.Range('E13').Value = '  -2.98%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '26.28'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -1.84%  '
$ws.Range('E15').Value = '  -3.38%  '
$ws.Range('D16').Value = '2.853.17'
$ws.Range('E16').Value = '  -0.87%  '
$ws.Range('D17').Value = '62.008.82'
$ws.Range('E17').Value = '  -0.73%  '
$ws.Range('D18').Value = '2.436.44'
$ws.Range('E18').Value = '  +0.27%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '10.81'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -3.68%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '7.11'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -2.20%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '327.50'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +0.29%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '4.08'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -2.29%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '1.98'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -3.78%  '
$ws.Range('E24').Value = '  -0.50%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '65.59'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +0.87%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '9.32'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +6.36%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '609.86'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -0.51%  '
$ws.Range('D28').Value = '2.565.85'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0942'
$ws.Range('E30').Value = '  -6.52%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '1.42'
$c.Style = "Normal"
$ws.Range('E31').Value = '  -5.28%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '7.95'
$c.Style = "Normal"
$ws.Range('E32').Value = '  -2.20%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.140'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +0.41%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.88'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -0.38%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '4.87'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -5.01%  '
$ws.Range('E36').Value = '  +0.16%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '1.42'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -5.71%  '
$ws.Range('E38').Value = '  -0.42%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '149.00'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +2.75%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '5.28'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -0.84%  '
$ws.Range('E41').Value = '  -2.32%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '1.72'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -2.78%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '42.51'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +1.56%  '
$ws.Range('E44').Value = '  +0.01%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '2.43'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -5.67%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '142.31'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -3.61%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '3.61'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -3.56%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '0.603'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +1.15%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.0521'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -1.47%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '19.37'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -7.15%  '
$ws.Range('D51').Value = '0.0₆0234'
$ws.Range('E51').Value = '  +8.74%  '
